# Add localized strings for the service-worker update modal
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KeyValuePairs")

# Keys first (column A for both new rows)
$ws.Range("A48").Value = "service_worker-update_headline"
$ws.Range("A49").Value = "service_worker-update_confirm_btn_txt"

# Row 48 translations
$ws.Range("B48").Value = "Update Available"
$ws.Range("C48").Value = "Update verfügbar"

# Row 49 translations
$ws.Range("B49").Value = "Update Now & Refresh"
$ws.Range("C49").Value = "Update installieren"

# Match the fill style used by the rest of the data rows (cellXfs index 1)
$ws.Range("A46:C46").Copy()
$ws.Range("A48:C49").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Expand the table to include the two new rows
$lo = $ws.ListObjects.Item("Tabelle2")
$lo.Resize($ws.Range("A1:C49"))

# Re-fit column A now that it holds a longer key string
$ws.Columns("A:A").ColumnWidth = 31.75

# Update the sheet/table view state to reflect scrolling to the new rows
[void]$ws.Range("A31").Select()
$ws.Application.ActiveWindow.ScrollRow = 31
[void]$ws.Range("D44").Select()

$wb.Save()
